$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.179.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.871.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.25"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.90"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.876.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "0.0000000"
$ws.Range("D14").Value = 0.0000253
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.531.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.890.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.261.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +9.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.112"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "486.84"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.746"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000170"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.029.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.80"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.35"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.55"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.844.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.142"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +11.34%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.323"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "430.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000277"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +20.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0364"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.87"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.26"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.70%  "
